$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Helper: write plain single-run text into a table cell ---
function Set-CellText($row, $colIndex, $text) {
    $cell = $row.Cells.Item($colIndex)
    $cell.Range.Text = $text
}

# --- Pins 4-7 (PWM rows): fill in the Color Sensor S0-S3 pins ---
# Pin 4 -> Color Sensor S0 / Red frequency reader
$row = $t.Rows.Item(6)
Set-CellText $row 3 "Color Sensor S0"
Set-CellText $row 4 "Red frequency reader"

# Pin 5 -> Color Sensor S1 / Green frequency reader
$row = $t.Rows.Item(7)
Set-CellText $row 3 "Color Sensor S1"
Set-CellText $row 4 "Green frequency reader"

# Pin 6 -> Color Sensor S2 (typed as two runs: "Color Sensor S" + "2") / Blue frequency reader
$row = $t.Rows.Item(8)
$cell = $row.Cells.Item(3)
$cell.Range.Text = "Color Sensor S"
$cell = $row.Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$tail = $d.Range($para.Range.Start, $para.Range.End - 1)
$tail.InsertAfter("2")
$tail.Bold = 1
$tail.Bold = 0
Set-CellText $row 4 "Blue frequency reader"

# Pin 7 -> Color Sensor S3 (typed as two runs: "Color Sensor S" + "3") / Clear frequency reader
$row = $t.Rows.Item(9)
$cell = $row.Cells.Item(3)
$cell.Range.Text = "Color Sensor S"
$cell = $row.Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$tail = $d.Range($para.Range.Start, $para.Range.End - 1)
$tail.InsertAfter("3")
$tail.Bold = 1
$tail.Bold = 0
Set-CellText $row 4 "Clear frequency reader"

# --- Clear the three "L298N Catch Arm" rows (pins 45, 47, 49) ---
$d.Content.Find.Execute("L298N Catch Arm", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Catch Arm ENA PWM", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Catch Arm IN1 Digital", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Catch Arm IN2 Digital", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null

# --- Clear the Color Sensor SCL / SDA rows (pins A4, A5) ---
$d.Content.Find.Execute("Color Sensor SCL", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Clock for Color Sensor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Color Sensor SDA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Data for Color Sensor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null

# --- Fill in the last row (pin A15): Distance Sensor ---
$row = $t.Rows.Item($t.Rows.Count)
Set-CellText $row 3 "Distance Sensor"
Set-CellText $row 4 "Reads analog distance sensor voltage"
